# Add a new "20-jul" data column (AE) to the right of the existing "19-jul"
# column (AD) on Sheet1, mirroring the header style and number formatting of
# the existing date columns, then fill in the new values for rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell: same style as AD1 (text-formatted header, centered via shared
# style index 1 in the original file -> numFmtId 49).
$ws.Range("AE1").Value = "20-jul"
$ws.Range("AE1").NumberFormat = $ws.Range("AD1").NumberFormat

# New column values for rows 2-11 (mirrors style of AD2:AD11 -> numFmtId 1,
# centered integer values).
$values = @(10, 16, 8, 13, 14, 14, 19, 15, 25, 23)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("AE$row")
    $cell.Value = $values[$i]
    $cell.NumberFormat = $ws.Range("AD$row").NumberFormat
    $cell.HorizontalAlignment = $ws.Range("AD$row").HorizontalAlignment
}

# Update the selection to match what Excel recorded after this edit.
$ws.Range("AH13").Select()
